# Update the "echo" TestCase sheet to add a ForEachCommand-style "pause"
# command (columns E:F) and change the "echo" command's second parameter
# column from a literal "admin" into a JSON target parameter, while the
# first parameter column becomes the "${name}" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2, column C: was blank, now holds the "${name}" placeholder.
# The existing format on C2 (Microsoft YaHei, color FF0451A5) already
# matches the target look, so only the value needs to change.
$ws.Range("C2").Value = "`${name}"

# --- New "pause" parameter cells D2:F2 and D3:F3 need a brand new font
# (MS PGothic, color FF0451A5). Copy the base formatting (borders,
# alignment, size) from C2 first, then override just the font.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2:F2").PasteSpecial(-4122) | Out-Null
$ws.Range("D3:F3").PasteSpecial(-4122) | Out-Null

# NOTE: a multi-area Range (comma-joined) only applies Font changes to its
# first area in this host, so the two rows are styled separately.
foreach ($rowRange in @($ws.Range("D2:F2"), $ws.Range("D3:F3"))) {
    $rowRange.Font.Name = "ＭＳ Ｐゴシック"
    $rowRange.Font.Size = 9
    $rowRange.Font.Color = 10834180
}

# --- Row 3, column C: used to hold "${name}"; now blank, restyled to
# match the "echo"/"pause" command-name header look (copy from D1).
$ws.Range("D1").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").ClearContents()

# --- Row 3, column D: used to hold the literal "admin"; now the JSON
# target parameter for the echo command.
$ws.Range("D3").Value = '{"target":"admin"}'

# --- Row 1: new "pause" command occupies columns E:F, mirroring "echo"
# in columns C:D. Copy the command-name header style from C1.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "pause"
$ws.Range("F1").Value = "pause"

# --- Row 2, column E: pause's first parameter, a literal number.
$ws.Range("E2").Value = 2000

# --- Row 3, column F: pause's second parameter, JSON target.
$ws.Range("F3").Value = '{"target":2000}'

# --- Column widths: column C keeps its old shared width; D:E get a new
# wider size for the pause columns, F narrower.
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 12.71
$ws.Range("F1").EntireColumn.ColumnWidth = 10.14

# --- Move the active selection to match the saved view state.
$ws.Range("G10").Select() | Out-Null

Write-Host "edit applied"
